$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a literal/text value into a cell without Excel's
# automatic number/date coercion, and without leaving behind a new
# number-format/quote-prefix style record (which a plain NumberFormat="@"
# or a leading "'" would do). We do this by writing a formula that
# evaluates to the exact text we want, then collapsing it to a plain
# value with a Copy / Paste-Special-Values round trip. ---
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# Numeric id / sort-order updates
$ws.Range("A2").Value = 55915885
$ws.Range("B2").Value = 103226

# Antal / Enhet / Ålder-Stadium (row previously had an empty "Antal" cell)
Set-TextValue $ws.Range("I2") "40"
$ws.Range("J2").Value = "stjälkar/strån/skott"
$ws.Range("K2").Value = "i frukt"
# Kön stays blank, but the cell itself becomes a present empty-text cell.
Set-TextValue $ws.Range("L2") ""

# Noggrannhet
$ws.Range("S2").Value = 10

# Externid removed entirely
$ws.Range("X2").ClearContents()

# Startdatum / Slutdatum - keep as literal text (not auto-converted to a date serial)
Set-TextValue $ws.Range("Y2") "2015-11-15"
Set-TextValue $ws.Range("AA2") "2015-11-15"

# Publik kommentar
$ws.Range("AC2").Value = "9 florala + 31 vegetativa skott. Denna lokal bör vara densamma som floraväktarlokalen ""C-Tie-0309 Gyllerboda, 450 m V om, Tierp s:n, Upl"", men nu med vettiga koordinater! Följearter: husmossa, väggmossa, lingon, vårfryle, blåbär."

# Biotop / Biotop-beskrivning (new)
$ws.Range("AH2").Value = "Lingontallskog"
$ws.Range("AI2").Value = "Talldominerad åsbarrskog"

# Rapportör / Observatörer / Projektnamn
$ws.Range("AW2").Value = "Sebastian Sundberg"
$ws.Range("AX2").Value = "Sebastian Sundberg, Owe Rosengren"
Set-TextValue $ws.Range("AY2") ""
